$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# Insert 16 new rows above row 5 (push existing rows 5-34 down to 21-50)
$ws.Rows("5:20").Insert()

$ws.Range("A5").Value = "Moving to the new World"
$ws.Range("A6").Value = "- sqlite - complete base functions to replace excel sheets"
$ws.Range("A7").Value = "- Update naming convensions"
$ws.Range("A8").Value = "- Move to new directory structure"
$ws.Range("A10").Value = "Incorporate Prod Date in App"
$ws.Range("A11").Value = "Incorporate Effective Date"
$ws.Range("A12").Value = "Add WIO concept"
$ws.Range("A13").Value = "Update GUI to look nice"
$ws.Range("A14").Value = "Update GUI to actually update"
$ws.Range("A15").Value = "More Calculations"
$ws.Range("A16").Value = "More Reports"
$ws.Range("A17").Value = "More Browses"
$ws.Range("A18").Value = "Mapping Stuff (Just show a little)"

$ws.Range("A6:A8").Style = "Normal"
$ws.Range("A6:A8").HorizontalAlignment = -4131
$ws.Range("A6:A8").WrapText = $true
$ws.Range("A6:A8").IndentLevel = 1

$ws.Range("B16").Select()
